# Model sheet: insert 3 blank rows above the "Revenue" row, pushing the
# Revenue/COGS/Gross Margin/R&D/SG&A/OpEx/Operating Income block (and the
# trailing spacer row) down from rows 3-11 to rows 6-14.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Model")
$ws.Rows("3:5").Insert()

# The edit also left the "Model" sheet active (vs. "Main" before), with the
# newly-inserted row 3 selected across the full row.
$ws.Activate()
[void]$ws.Rows(3).Select()
